# Apply the crypto price/volume refresh produced by the scheduled GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.481.72'
$ws.Range("E2").Value = '  +2.00%  '
$ws.Range("D3").Value = '1.865.25'
$ws.Range("E3").Value = '  +1.02%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '''311.64'
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("D7").Value = '''0.4771'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '''0.3740'
$ws.Range("E8").Value = '  +2.06%  '
$ws.Range("D9").Value = '''0.07319'
$ws.Range("E9").Value = '  +1.33%  '
$ws.Range("D10").Value = '''0.9341'
$ws.Range("E10").Value = '  +0.74%  '
$ws.Range("D11").Value = '''20.65'
$ws.Range("E11").Value = '  +4.90%  '
$ws.Range("D12").Value = '''0.07823'
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("D13").Value = '1.859.11'
$ws.Range("E13").Value = '  -1.54%  '
$ws.Range("D14").Value = '''5.435'
$ws.Range("E14").Value = '  +2.24%  '
$ws.Range("D15").Value = '''6.547'
$ws.Range("E15").Value = '  +2.13%  '
$ws.Range("D16").Value = '''90.22'
$ws.Range("E16").Value = '  +1.59%  '
$ws.Range("D17").Value = '''1.013'
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").Value = '''0.000008881'
$ws.Range("E18").Value = '  +2.83%  '
$ws.Range("D19").Value = '''1.011'
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").Value = '27.533.91'
$ws.Range("E20").Value = '  +2.09%  '
$ws.Range("D21").Value = '''14.56'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '''5.111'
$ws.Range("E22").Value = '  +1.10%  '
$ws.Range("D23").Value = '''10.69'
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("D24").Value = '''1.945'
$ws.Range("E24").Value = '  +0.86%  '
$ws.Range("D25").Value = '''154.71'
$ws.Range("E25").Value = '  +1.52%  '
$ws.Range("D26").Value = '''18.45'
$ws.Range("E26").Value = '  +1.49%  '
$ws.Range("D27").Value = '''2.020'
$ws.Range("E27").Value = '  +1.29%  '
$ws.Range("D28").Value = '''115.50'
$ws.Range("E28").Value = '  +1.18%  '
$ws.Range("D29").Value = '''4.979'
$ws.Range("E29").Value = '  +0.44%  '
$ws.Range("D30").Value = '''0.08895'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").Value = '''3.334'
$ws.Range("E31").Value = '  +0.44%  '
$ws.Range("D32").Value = '''1.218'
$ws.Range("E32").Value = '  +3.98%  '
$ws.Range("D33").Value = '''0.7577'
$ws.Range("E33").Value = '  +1.69%  '
$ws.Range("D34").Value = '''4.609'
$ws.Range("E34").Value = '  +2.62%  '
$ws.Range("D35").Value = '''2.738'
$ws.Range("E35").Value = '  +0.56%  '
$ws.Range("D36").Value = '''0.02041'
$ws.Range("E36").Value = '  +4.33%  '
$ws.Range("D37").Value = '''1.120'
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("D38").Value = '''2.995'
$ws.Range("E38").Value = '  +0.29%  '
$ws.Range("D39").Value = '''0.05262'
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("D40").Value = '''0.5304'
$ws.Range("E40").Value = '  +2.14%  '
$ws.Range("D41").Value = '''7.070'
$ws.Range("E41").Value = '  +1.13%  '
$ws.Range("D42").Value = '''0.1522'
$ws.Range("E42").Value = '  +0.82%  '
$ws.Range("D43").Value = '''8.461'
$ws.Range("E43").Value = '  +3.19%  '
$ws.Range("D44").Value = '''10.66'
$ws.Range("E44").Value = '  +0.53%  '
$ws.Range("D45").Value = '''0.4798'
$ws.Range("E45").Value = '  +1.45%  '
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '''102.78'
$ws.Range("E47").Value = '  +1.23%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '''1.650'
$ws.Range("E48").Value = '  +2.99%  '
$ws.Range("D49").Value = '''67.36'
$ws.Range("E49").Value = '  +2.96%  '
$ws.Range("D51").Value = '''0.9174'
$ws.Range("E51").Value = '  +3.49%  '
